$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-14 20:18:54"
$ws.Range("I2").Value = "35.0 mm"
$ws.Range("E3").Value = "2026-02-14 20:18:57"
$ws.Range("I3").Value = "14.9 mm"
$ws.Range("N3").Value = "-6.9 °C 19:59 TU"
$ws.Range("E4").Value = "2026-02-14 20:19:00"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "73%"
$ws.Range("J4").Value = "996.8 hPa"
$ws.Range("N4").Value = "5.2 °C 19:52 TU"
$ws.Range("O4").Value = "10.8 °C"
$ws.Range("E5").Value = "2026-02-14 20:19:03"
$ws.Range("I5").Value = "21.6 mm"
$ws.Range("N5").Value = "-6.8 °C 19:56 TU"
$ws.Range("E6").Value = "2026-02-14 20:19:05"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "76%"
$ws.Range("J6").Value = "996.8 hPa"
$ws.Range("E7").Value = "2026-02-14 20:19:08"
$ws.Range("J7").Value = "997.0 hPa"
$ws.Range("O7").Value = "13.3 °C"
$ws.Range("E8").Value = "2026-02-14 20:19:11"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "63%"
$ws.Range("J8").Value = "996.8 hPa"
$ws.Range("E9").Value = "2026-02-14 20:19:13"
$ws.Range("E10").Value = "2026-02-14 20:19:16"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "77%"
$ws.Range("E11").Value = "2026-02-14 20:19:18"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "60%"
$ws.Range("E12").Value = "2026-02-14 20:19:21"
$ws.Range("E13").Value = "2026-02-14 20:19:24"
$ws.Range("J13").Value = "999.6 hPa"
$ws.Range("E14").Value = "2026-02-14 20:19:27"
$ws.Range("E15").Value = "2026-02-14 20:19:30"
$ws.Range("N15").Value = "10.0 °C 19:59 TU"
$ws.Range("O15").Value = "11.3 °C"
$ws.Range("E16").Value = "2026-02-14 20:19:32"
$ws.Range("O16").Value = "-6.0 °C"
$ws.Range("E17").Value = "2026-02-14 20:19:35"
$ws.Range("O17").Value = "1.8 °C"
$ws.Range("E18").Value = "2026-02-14 20:19:38"
$ws.Range("J18").Value = "997.0 hPa"
$ws.Range("O18").Value = "10.6 °C"
$ws.Range("E19").Value = "2026-02-14 20:19:41"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "78%"
$ws.Range("E20").Value = "2026-02-14 20:19:43"
$ws.Range("I20").Value = "4.1 mm"
$ws.Range("N20").Value = "-7.6 °C 19:59 TU"
$ws.Range("O20").Value = "-5.3 °C"
$ws.Range("E21").Value = "2026-02-14 20:19:46"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "69%"
$ws.Range("J21").Value = "999.5 hPa"
$ws.Range("E22").Value = "2026-02-14 20:19:49"
$ws.Range("N22").Value = "-9.1 °C 19:56 TU"
$ws.Range("O22").Value = "-6.8 °C"
$ws.Range("E23").Value = "2026-02-14 20:19:51"
$ws.Range("I23").Value = "39.1 mm"
$ws.Range("O23").Value = "-6.0 °C"
$ws.Range("E24").Value = "2026-02-14 20:19:54"
$ws.Range("J24").Value = "1001.1 hPa"
$ws.Range("E25").Value = "2026-02-14 20:19:57"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "84%"
$ws.Range("I25").Value = "14.4 mm"
$ws.Range("N25").Value = "-7.2 °C 19:59 TU"
$ws.Range("E26").Value = "2026-02-14 20:20:00"
$ws.Range("E27").Value = "2026-02-14 20:20:03"
$ws.Range("N27").Value = "-6.2 °C 19:55 TU"
$ws.Range("O27").Value = "-3.0 °C"
$ws.Range("E28").Value = "2026-02-14 20:20:05"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "68%"
$ws.Range("J28").Value = "996.8 hPa"
$ws.Range("E29").Value = "2026-02-14 20:20:08"
$ws.Range("E30").Value = "2026-02-14 20:20:11"
$ws.Range("J30").Value = "996.7 hPa"
$ws.Range("L30").Value = "96.5 km/h - 8º 19:57 TU"
$ws.Range("E31").Value = "2026-02-14 20:20:14"
$ws.Range("J31").Value = "995.8 hPa"
$ws.Range("N31").Value = "8.5 °C 19:43 TU"
$ws.Range("E32").Value = "2026-02-14 20:20:17"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "84%"
$ws.Range("E33").Value = "2026-02-14 20:20:19"
$ws.Range("J33").Value = "998.9 hPa"
$ws.Range("O33").Value = "4.0 °C"
$ws.Range("E34").Value = "2026-02-14 20:20:22"
$ws.Range("I34").Value = "3.2 mm"
$ws.Range("N34").Value = "-4.8 °C 19:59 TU"
$ws.Range("E35").Value = "2026-02-14 20:20:25"
$ws.Range("J35").Value = "1003.8 hPa"
$ws.Range("O35").Value = "2.7 °C"
$ws.Range("E36").Value = "2026-02-14 20:20:28"
$ws.Range("J36").Value = "997.5 hPa"
$ws.Range("N36").Value = "10.7 °C 19:59 TU"
$ws.Range("O36").Value = "12.0 °C"
$ws.Range("E37").Value = "2026-02-14 20:20:31"
$ws.Range("J37").Value = "997.8 hPa"
$ws.Range("E38").Value = "2026-02-14 20:20:34"
$ws.Range("N38").Value = "6.3 °C 19:55 TU"
$ws.Range("O38").Value = "10.2 °C"
$ws.Range("E39").Value = "2026-02-14 20:20:36"
$ws.Range("I39").Value = "12.5 mm"
$ws.Range("N39").Value = "-8.5 °C 19:57 TU"
$ws.Range("O39").Value = "-5.7 °C"
$ws.Range("E40").Value = "2026-02-14 20:20:39"
$ws.Range("I40").Value = "0.7 mm"
$ws.Range("J40").Value = "1000.1 hPa"
$ws.Range("E41").Value = "2026-02-14 20:20:42"
$ws.Range("J41").Value = "998.8 hPa"
$ws.Range("E42").Value = "2026-02-14 20:20:45"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "62%"
$ws.Range("E43").Value = "2026-02-14 20:20:48"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "67%"
$ws.Range("O43").Value = "9.2 °C"
$ws.Range("E44").Value = "2026-02-14 20:20:51"
$ws.Range("G44").Value = "271 cm"
$ws.Range("I44").Value = "36.8 mm"
$ws.Range("N44").Value = "-7.0 °C 19:56 TU"
$ws.Range("E45").Value = "2026-02-14 20:20:53"
$ws.Range("J45").Value = "1006.0 hPa"
$ws.Range("E46").Value = "2026-02-14 20:20:56"
